$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New xG_home (D), xG_away (E), goals_home (F), goals_away (G) results
# for the matches in rows 10-15 (previously blank in those columns).
$data = @{
    10 = @(1.44548, 2.7383, 2, 2)
    11 = @(1.7297, 0.437896, 1, 2)
    12 = @(1.69279, 0.824038, 4, 1)
    13 = @(2.66532, 1.61165, 2, 2)
    14 = @(1.26331, 0.779401, 2, 1)
    15 = @(0.530442, 1.24789, 1, 2)
}

$rows = 10..15

# Fill column-by-column (D, then E, then F, then G) to mirror how the
# source data was appended to the sheet.
foreach ($row in $rows) {
    $c = $ws.Range("D$row")
    $c.Value = "'" + $data[$row][0]
    $c.Style = "Normal"
}
foreach ($row in $rows) {
    $c = $ws.Range("E$row")
    $c.Value = "'" + $data[$row][1]
    $c.Style = "Normal"
}
foreach ($row in $rows) {
    $c = $ws.Range("F$row")
    $c.Value = "'" + $data[$row][2]
    $c.Style = "Normal"
}
foreach ($row in $rows) {
    $c = $ws.Range("G$row")
    $c.Value = "'" + $data[$row][3]
    $c.Style = "Normal"
}
